$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link) ---
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B12").Value = "MCDex"
$ws.Range("C12").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

# --- Numeric-looking text columns (Price / Volume) need a leading apostrophe
#     so Excel stores them as text, matching the source workbook formatting ---
$ws.Range("D2").Value = "'329.05"
$ws.Range("D3").Value = "'43.91"
$ws.Range("E3").Value = "'-1.52%"
$ws.Range("D4").Value = "'5.469"
$ws.Range("E4").Value = "'-1.92%"
$ws.Range("D5").Value = "'0.07981"
$ws.Range("E5").Value = "'-1.19%"
$ws.Range("D6").Value = "'1.982"
$ws.Range("E6").Value = "'4.07%"
$ws.Range("D7").Value = "'4.376"
$ws.Range("E7").Value = "'1.94%"
$ws.Range("D8").Value = "'2.574"
$ws.Range("E8").Value = "'-4.32%"
$ws.Range("D9").Value = "'0.9498"
$ws.Range("E9").Value = "'0.77%"
$ws.Range("D10").Value = "'0.1110"
$ws.Range("E10").Value = "'-5.22%"
$ws.Range("D11").Value = "'0.1888"
$ws.Range("E11").Value = "'1.70%"
$ws.Range("D12").Value = "'10.61"
$ws.Range("E12").Value = "'25.71%"
$ws.Range("D13").Value = "'0.09951"
$ws.Range("E13").Value = "'-0.51%"
$ws.Range("D14").Value = "'0.04845"
$ws.Range("E14").Value = "'13.25%"
$ws.Range("D15").Value = "'0.1064"
$ws.Range("E15").Value = "'-0.17%"
$ws.Range("D16").Value = "'0.001272"
$ws.Range("E16").Value = "'-0.59%"
$ws.Range("D17").Value = "'0.04088"
$ws.Range("E17").Value = "'-2.49%"
$ws.Range("D18").Value = "'0.005961"
$ws.Range("E18").Value = "'1.35%"
$ws.Range("D19").Value = "'3.368"
$ws.Range("E19").Value = "'-5.82%"
$ws.Range("D20").Value = "'0.3466"
$ws.Range("E20").Value = "'-0.93%"
$ws.Range("D21").Value = "'0.1420"
$ws.Range("E21").Value = "'3.71%"
$ws.Range("E23").Value = "'2.10%"
$ws.Range("D24").Value = "'0.004365"
$ws.Range("E24").Value = "'-3.23%"
$ws.Range("E25").Value = "'1.60%"
$ws.Range("D26").Value = "'0.0003747"
$ws.Range("E26").Value = "'-6.07%"
$ws.Range("D38").Value = "'0.02569"
$ws.Range("E38").Value = "'-2.65%"
$ws.Range("D39").Value = "'0.05642"
$ws.Range("E39").Value = "'3.64%"
$ws.Range("D40").Value = "'0.007563"
$ws.Range("E40").Value = "'-1.19%"
$ws.Range("E41").Value = "'-0.10%"
$ws.Range("D42").Value = "'0.007393"
$ws.Range("E42").Value = "'4.57%"
$ws.Range("D43").Value = "'0.002016"
$ws.Range("E43").Value = "'-0.50%"
$ws.Range("D44").Value = "'0.008609"
$ws.Range("E44").Value = "'-2.73%"
$ws.Range("D45").Value = "'0.00007133"
$ws.Range("E45").Value = "'-0.49%"
$ws.Range("E46").Value = "'-0.06%"
$ws.Range("D47").Value = "'0.003532"
$ws.Range("E47").Value = "'55.62%"
$ws.Range("D48").Value = "'0.003725"
$ws.Range("E48").Value = "'1.66%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'-0.06%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'-0.06%"
